# Slide 5 contains the coding-security table. In row 2 (Cross Site
# Scripting row) of that table:
#  - column 4 ("적용여부 (O/X)") changes from " X" to " O"
#  - column 5 ("비고") gets two new paragraphs of explanatory text
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item(2)
$tbl = $sh.Table

# Column 4: " X" -> " O"
$cellApply = $tbl.Cell(2, 4)
$trApply = $cellApply.Shape.TextFrame.TextRange
$trApply.Text = " O"

# Column 5: add remark text (two paragraphs)
$cellRemark = $tbl.Cell(2, 5)
$trRemark = $cellRemark.Shape.TextFrame.TextRange
$trRemark.Text = "게시판,자료실,웹메일 미사용`r검색창에 < > ( ) # & 입력불가처리"
